$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.382.80"
$ws.Range("E2").Value = "  -4.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.701.76"
$ws.Range("E3").Value = "  -4.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.89"
$ws.Range("E5").Value = "  -2.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.46"
$ws.Range("E6").Value = "  +7.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.697.32"
$ws.Range("E7").Value = "  -4.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.628"
$ws.Range("E8").Value = "  -7.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.717"
$ws.Range("E10").Value = "  -5.76%  "
$ws.Range("E11").Value = "  -8.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.82"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("E13").Value = "  -9.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.41"
$ws.Range("E14").Value = "  -7.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.211.96"
$ws.Range("E15").Value = "  -6.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.717.13"
$ws.Range("E16").Value = "  -4.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.49"
$ws.Range("E17").Value = "  -6.79%  "
$ws.Range("E18").Value = "  -2.64%  "
$ws.Range("B19").Value = "Polygon"
$ws.Range("C19").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.12"
$ws.Range("E19").Value = "  -7.58%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.76"
$ws.Range("E20").Value = "  -8.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.258.35"
$ws.Range("E21").Value = "  -4.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "409.05"
$ws.Range("E22").Value = "  -6.33%  "
$ws.Range("E23").Value = "  -5.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.57"
$ws.Range("E24").Value = "  -6.56%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.19"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.81"
$ws.Range("E27").Value = "  -7.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.92"
$ws.Range("E28").Value = "  -4.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.07"
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.48"
$ws.Range("E30").Value = "  -7.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.62"
$ws.Range("E31").Value = "  -7.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.45"
$ws.Range("E32").Value = "  -6.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.55"
$ws.Range("E33").Value = "  -7.73%  "
$ws.Range("E34").Value = "  -7.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "65.37"
$ws.Range("E35").Value = "  -5.12%  "
$ws.Range("E36").Value = "  -16.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "601.75"
$ws.Range("E37").Value = "  -3.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0897"
$ws.Range("E38").Value = "  -9.24%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.401"
$ws.Range("E40").Value = "  -5.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.137"
$ws.Range("E42").Value = "  -4.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.77"
$ws.Range("E43").Value = "  +2.99%  "
$ws.Range("E44").Value = "  -9.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.94"
$ws.Range("E45").Value = "  -9.59%  "
$ws.Range("E46").Value = "  -7.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.29"
$ws.Range("E47").Value = "  -9.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.791.20"
$ws.Range("E48").Value = "  -1.81%  "
$ws.Range("E49").Value = "  -7.70%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.17"
$ws.Range("E50").Value = "  -6.69%  "
$ws.Range("B51").Value = "WEMIXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.69"
$ws.Range("E51").Value = "  -4.00%  "
